$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition list) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7438
$ws1.Range("F3").Value = 68
$ws1.Range("F4").Value = 211
$ws1.Range("F5").Value = 218
$ws1.Range("F6").Value = 1119
$ws1.Range("F7").Value = 195
$ws1.Range("F8").Value = 14
$ws1.Range("F9").Value = 115

# --- Sheet "全部类型" (all types combined list) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7438
$ws4.Range("F3").Value = 68
$ws4.Range("F4").Value = 211
$ws4.Range("F5").Value = 218
$ws4.Range("F6").Value = 1119
$ws4.Range("F7").Value = 195
$ws4.Range("F9").Value = 14
$ws4.Range("F10").Value = 115
